# Populate the "Transport Order List" distance matrix on the
# "Transport Order Data" sheet (xl/worksheets/sheet3.xml).
#
# This reproduces a lower-triangular distance matrix in columns O:T,
# rows 10-14, using the new location/supplier labels:
#   Nuremburg, Munich, Stuttgart, Supplier S.A (Porto), Supplier Lda (Barcelona)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transport Order Data")

# Row 10 - Nuremburg
$ws.Cells.Item(10, 15).Value = "Nuremburg"
$ws.Cells.Item(10, 16).Value = 0

# Row 11 - Munich
$ws.Cells.Item(11, 15).Value = "Munich"
$ws.Cells.Item(11, 16).Value = 170
$ws.Cells.Item(11, 17).Value = 0

# Row 12 - Stuttgart
$ws.Cells.Item(12, 15).Value = "Stuttgart"
$ws.Cells.Item(12, 16).Value = 210
$ws.Cells.Item(12, 17).Value = 243
$ws.Cells.Item(12, 18).Value = 0

# Row 13 - Supplier S.A (Porto)
$ws.Cells.Item(13, 15).Value = "Supplier S.A (Porto)"
$ws.Cells.Item(13, 16).Value = 2219
$ws.Cells.Item(13, 17).Value = 2253
$ws.Cells.Item(13, 18).Value = 2042
$ws.Cells.Item(13, 19).Value = 0

# Row 14 - Supplier Lda (Barcelona)
$ws.Cells.Item(14, 15).Value = "Supplier Lda (Barcelona)"
$ws.Cells.Item(14, 16).Value = 1444
$ws.Cells.Item(14, 17).Value = 1369
$ws.Cells.Item(14, 18).Value = 1267
$ws.Cells.Item(14, 19).Value = 1127
$ws.Cells.Item(14, 20).Value = 0

# Update the selection to reflect where the user ended up working
$ws.Activate() | Out-Null
$ws.Range("P15").Select() | Out-Null
